$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.940.54'
$ws.Range("E2").Value = '  -0.64%  '
$ws.Range("D3").Value = '2.213.76'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '73.30'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.84%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0955'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.13'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").Value = '2.545.08'
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.842'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.00%  '
$ws.Range("D17").Value = '2.194.85'
$ws.Range("E17").Value = '  -3.53%  '
$ws.Range("D18").Value = '41.895.66'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("E19").Value = '  +10.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.15'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("E22").Value = '  +21.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.94%  '
$ws.Range("E24").Value = '  -7.28%  '
$ws.Range("E25").Value = '  +0.06%  '
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("E28").Value = '  -1.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.91%  '
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("E32").Value = '  +8.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0793'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.46%  '
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '28.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.39%  '
$ws.Range("E36").Value = '  -7.85%  '
$ws.Range("E37").Value = '  -5.20%  '
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.84'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '66.18'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.52%  '
$ws.Range("E41").Value = '  -3.14%  '
$ws.Range("E42").Value = '  -2.80%  '
$ws.Range("E43").Value = '  -0.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.100'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.13%  '
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").Value = '2.421.69'
$ws.Range("E51").Value = '  -1.23%  '
